$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attendance fill-in data for weeks 7-10 (columns M,N,O,P) across the 63 student rows (7-69).
# P = Present, A = Absent. Values were previously blank for these columns.
$data = @(
    @{ Row = 7; M = "P"; N = "P"; O = "A"; P = "P" },
    @{ Row = 8; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 9; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 10; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 11; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 12; M = "A"; N = "A"; O = "P"; P = "P" },
    @{ Row = 13; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 14; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 15; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 16; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 17; M = "A"; N = "A"; O = "P"; P = "P" },
    @{ Row = 18; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 19; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 20; M = "P"; N = "A"; O = "P"; P = "A" },
    @{ Row = 21; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 22; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 23; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 24; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 25; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 26; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 27; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 28; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 29; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 30; M = "P"; N = "A"; O = "P"; P = "A" },
    @{ Row = 31; M = "P"; N = "P"; O = "A"; P = "A" },
    @{ Row = 32; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 33; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 34; M = "P"; N = "P"; O = "A"; P = "P" },
    @{ Row = 35; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 36; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 37; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 38; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 39; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 40; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 41; M = "A"; N = "A"; O = "P"; P = "P" },
    @{ Row = 42; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 43; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 44; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 45; M = "A"; N = "A"; O = "A"; P = "A" },
    @{ Row = 46; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 47; M = "A"; N = "P"; O = "A"; P = "A" },
    @{ Row = 48; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 49; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 50; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 51; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 52; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 53; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 54; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 55; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 56; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 57; M = "P"; N = "A"; O = "P"; P = "A" },
    @{ Row = 58; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 59; M = "A"; N = "P"; O = "P"; P = "A" },
    @{ Row = 60; M = "A"; N = "P"; O = "A"; P = "P" },
    @{ Row = 61; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 62; M = "A"; N = "P"; O = "P"; P = "P" },
    @{ Row = 63; M = "P"; N = "A"; O = "P"; P = "P" },
    @{ Row = 64; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 65; M = "P"; N = "P"; O = "P"; P = "P" },
    @{ Row = 66; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 67; M = "P"; N = "P"; O = "P"; P = "A" },
    @{ Row = 68; M = "A"; N = "A"; O = "P"; P = "P" },
    @{ Row = 69; M = "P"; N = "P"; O = "P"; P = "P" }
)

foreach ($entry in $data) {
    $r = $entry.Row

    # Copy the cell format from column L (already carries the attendance-grid
    # border/alignment styling) onto M:P before writing the new marks, so the
    # newly-populated cells pick up the same "data" style as the rest of the row
    # instead of keeping the blank-placeholder style.
    $ws.Range("L$r").Copy()
    $ws.Range("M$r`:P$r").PasteSpecial(-4122)

    $ws.Range("M$r").Value = $entry.M
    $ws.Range("N$r").Value = $entry.N
    $ws.Range("O$r").Value = $entry.O
    $ws.Range("P$r").Value = $entry.P
}

# Widen column D slightly (Remark column) to fit content.
$ws.Columns("D").ColumnWidth = 11.25
